$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "Area" / "Atotal" columns, plus a small J:K summary block ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Column D: re-derive the shared "midpoint depth" formula across D3:D8 ---
$ws.Range("D3:D8").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- Column E: re-apply formula (text unchanged, but now shares a (new) group) ---
$ws.Range("E4:E8").Formula = "=(D4-D3)*(B4/100)*C4"

# --- Column G: incremental area per segment ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Column H: total area ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Columns J/K: quick-reference copies of the totals ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Final selection, matching the saved UI state ---
$ws.Range("J2:K2").Select() | Out-Null
